$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 8099.6665
$ws.Range("I64").Value = 8399.799999999999
$ws.Range("K64").Value = 8399.799999999999
$ws.Range("M64").Value = -8151.799999999999
$ws.Range("H67").Value = 8099.6665
$ws.Range("I67").Value = 8399.799999999999
$ws.Range("K67").Value = 8399.799999999999
$ws.Range("M67").Value = -7541.799999999999
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H138").Value = 24392692
$ws.Range("I138").Value = 1132.4445
$ws.Range("J138").Value = 43481740
$ws.Range("K138").Value = 3397.3335
$ws.Range("L138").Value = 130445220
$ws.Range("M138").Value = 1742.6665
$ws.Range("N138").Value = -130455500

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5992.8335
$ws.Range("J2").Value = 7983
$ws.Range("L2").Value = 7983
$ws.Range("N2").Value = -8209
$ws.Range("H25").Value = 2906.2856
$ws.Range("I25").Value = 211
$ws.Range("J25").Value = 6500
$ws.Range("K25").Value = 211
$ws.Range("L25").Value = 6500
$ws.Range("M25").Value = 191
$ws.Range("N25").Value = -7304
$ws.Range("H39").Value = 8000
$ws.Range("I39").Value = 9500
$ws.Range("K39").Value = 9500
$ws.Range("M39").Value = -8980
$ws.Range("H45").Value = 3150
$ws.Range("I45").Value = 2528.5715
$ws.Range("K45").Value = 2528.5715
$ws.Range("M45").Value = -2151.5715
$ws.Range("H74").Value = 5538.613
$ws.Range("I74").Value = 5368.8623
$ws.Range("K74").Value = 5368.8623
$ws.Range("M74").Value = -4494.8623
$ws.Range("H77").Value = 5538.613
$ws.Range("I77").Value = 5368.8623
$ws.Range("K77").Value = 26844.3115
$ws.Range("M77").Value = -22476.3115
$ws.Range("H116").Value = 5992.8335
$ws.Range("J116").Value = 7983
$ws.Range("L116").Value = 7983
$ws.Range("N116").Value = -12571
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H124").Value = 24032.428
$ws.Range("J124").Value = 24032.428
$ws.Range("L124").Value = 24032.428
$ws.Range("N124").Value = -33852.428
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840
$ws.Range("H135").Value = 109338.695
$ws.Range("J135").Value = 109338.695
$ws.Range("L135").Value = 109338.695
$ws.Range("N135").Value = -119478.695

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5992.8335
$ws.Range("J3").Value = 7983
$ws.Range("L3").Value = 7983
$ws.Range("N3").Value = -8211
$ws.Range("H86").Value = 14109.647
$ws.Range("I86").Value = 2252.3
$ws.Range("K86").Value = 2252.3
$ws.Range("M86").Value = -1129.3
$ws.Range("H89").Value = 14109.647
$ws.Range("I89").Value = 2252.3
$ws.Range("K89").Value = 11261.5
$ws.Range("M89").Value = -5645.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H134").Value = 944.5714
$ws.Range("I134").Value = 940.3077
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 2820.9231
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -285.9231
$ws.Range("N134").Value = -8070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2531.4375
$ws.Range("I58").Value = 2160.25
$ws.Range("K58").Value = 2160.25
$ws.Range("M58").Value = -1957.25
$ws.Range("H86").Value = 10641.286
$ws.Range("I86").Value = 8663.333000000001
$ws.Range("J86").Value = 12124.75
$ws.Range("K86").Value = 8663.333000000001
$ws.Range("L86").Value = 12124.75
$ws.Range("M86").Value = -7540.333000000001
$ws.Range("N86").Value = -14370.75
$ws.Range("H89").Value = 10641.286
$ws.Range("I89").Value = 8663.333000000001
$ws.Range("J89").Value = 12124.75
$ws.Range("K89").Value = 43316.665
$ws.Range("L89").Value = 60623.75
$ws.Range("M89").Value = -37700.665
$ws.Range("N89").Value = -71855.75
$ws.Range("H132").Value = 4668.6665
$ws.Range("I132").Value = 4495.4736
$ws.Range("J132").Value = 5080
$ws.Range("K132").Value = 13486.4208
$ws.Range("L132").Value = 15240
$ws.Range("M132").Value = -10956.4208
$ws.Range("N132").Value = -20300
$ws.Range("H134").Value = 3328.55
$ws.Range("I134").Value = 2926.923
$ws.Range("K134").Value = 8780.769
$ws.Range("M134").Value = -6245.769
$ws.Range("H136").Value = 2531.4375
$ws.Range("I136").Value = 2160.25
$ws.Range("K136").Value = 6480.75
$ws.Range("M136").Value = -3930.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1668081.1
$ws.Range("J81").Value = 1749.25
$ws.Range("L81").Value = 5247.75
$ws.Range("N81").Value = -7493.75
$ws.Range("H84").Value = 1668081.1
$ws.Range("J84").Value = 1749.25
$ws.Range("L84").Value = 15743.25
$ws.Range("N84").Value = -26975.25
$ws.Range("H133").Value = 1593.6666
$ws.Range("I133").Value = 882.25
$ws.Range("K133").Value = 2646.75
$ws.Range("M133").Value = 2413.25
$ws.Range("H137").Value = 5073.636
$ws.Range("J137").Value = 5746.5
$ws.Range("L137").Value = 17239.5
$ws.Range("N137").Value = -27439.5
$ws.Range("H141").Value = 9155.714
$ws.Range("I141").Value = 1686
$ws.Range("K141").Value = 5058
$ws.Range("M141").Value = 122

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2667.2
$ws.Range("I132").Value = 2389.8333
$ws.Range("K132").Value = 7169.499899999999
$ws.Range("M132").Value = -4639.499899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4765.727
$ws.Range("I7").Value = 3996.7144
$ws.Range("J7").Value = 6111.5
$ws.Range("K7").Value = 3996.7144
$ws.Range("L7").Value = 6111.5
$ws.Range("M7").Value = -3884.7144
$ws.Range("N7").Value = -6335.5
$ws.Range("H123").Value = 55999.5
$ws.Range("J123").Value = 74999
$ws.Range("L123").Value = 74999
$ws.Range("N123").Value = -84799
$ws.Range("H126").Value = 4765.727
$ws.Range("I126").Value = 3996.7144
$ws.Range("J126").Value = 6111.5
$ws.Range("K126").Value = 11990.1432
$ws.Range("L126").Value = 18334.5
$ws.Range("M126").Value = -9520.143199999999
$ws.Range("N126").Value = -23274.5
$ws.Range("H136").Value = 3439.6667
$ws.Range("I136").Value = 3848.3333
$ws.Range("J136").Value = 1396.3334
$ws.Range("K136").Value = 11544.9999
$ws.Range("L136").Value = 4189.0002
$ws.Range("M136").Value = -8994.999899999999
$ws.Range("N136").Value = -9289.0002

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 977.5714
$ws.Range("I132").Value = 832.2727
$ws.Range("J132").Value = 1510.3334
$ws.Range("K132").Value = 2496.8181
$ws.Range("L132").Value = 4531.0002
$ws.Range("M132").Value = 33.18190000000004
$ws.Range("N132").Value = -9591.0002
